$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6, shifting existing rows 6-33 down to 7-34
$ws.Rows.Item(6).Insert()

# Fill in the new row's data (Applicateur d'encolure)
$ws.Range("A6").Value = 423700
$ws.Range("B6").Value = "Applicateur d'encolure"
$ws.Range("C6").Value = 2190
$ws.Range("D6").Value = 2190
$ws.Range("E6").Value = 1642.5
$ws.Range("F6").Value = 775

# Copy formatting from row 7 into row 6 (styles for C/D/E = style 1, F = style 3)
$ws.Range("C7:F7").Copy()
$ws.Range("C6:F6").PasteSpecial(-4122)  # xlPasteFormats

# F1 no longer carries the extra number-format style (reverts to default/general)
$ws.Range("F1").Style = "Normal"

# Reset the selection to E6 as in the target sheetView
$ws.Range("E6").Select()
